$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D header
$ws.Range("D1").Value = "ITI"

# Updated data values (Trial, Question, ConditionType, ITI) for rows 2-17 (Trial 1-16)
$data = @(
    @(1, 7, 1, 6),
    @(2, 31, 2, 8),
    @(3, 20, 2, 6),
    @(4, 17, 1, 8),
    @(5, 26, 4, 7),
    @(6, 32, 2, 6),
    @(7, 37, 4, 8),
    @(8, 6, 4, 8),
    @(9, 38, 4, 6),
    @(10, 19, 3, 8),
    @(11, 21, 3, 6),
    @(12, 5, 2, 6),
    @(13, 2, 3, 8),
    @(14, 28, 1, 8),
    @(15, 3, 1, 7),
    @(16, 12, 3, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

# Remove the now-unused trailing rows (previously Trial 17-19, rows 18-20)
$ws.Range("A18:D20").Delete()

$ws.Range("F8").Select()
